# Auto-generated edit script: update Chocobo_Profits value columns (H-N) per scraped diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 2300.2
$ws.Range("I28").Value = 2625.25
$ws.Range("J28").Value = 1000
$ws.Range("K28").Value = 2625.25
$ws.Range("L28").Value = 1000
$ws.Range("M28").Value = -2140.25
$ws.Range("N28").Value = -1970
$ws.Range("H112").Value = 1281.7377
$ws.Range("J112").Value = 1281.7377
$ws.Range("L112").Value = 3845.2131
$ws.Range("N112").Value = -6061.2131
$ws.Range("H132").Value = 28059874
$ws.Range("I132").Value = 29415896
$ws.Range("K132").Value = 88247688
$ws.Range("M132").Value = -88245158
$ws.Range("H137").Value = 3971140.8
$ws.Range("I137").Value = 5291854.5
$ws.Range("J137").Value = 9000
$ws.Range("K137").Value = 15875563.5
$ws.Range("L137").Value = 27000
$ws.Range("M137").Value = -15873013.5
$ws.Range("N137").Value = -32100

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 513.75
$ws.Range("I2").Value = 513.75
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 513.75
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -400.75
$ws.Range("N2").ClearContents()
$ws.Range("H45").Value = 4002.6667
$ws.Range("I45").Value = 4403.2
$ws.Range("J45").Value = 2000
$ws.Range("K45").Value = 4403.2
$ws.Range("L45").Value = 2000
$ws.Range("M45").Value = -4026.2
$ws.Range("N45").Value = -2754
$ws.Range("H74").Value = 2362
$ws.Range("I74").Value = 1470.4286
$ws.Range("K74").Value = 1470.4286
$ws.Range("M74").Value = -596.4286
$ws.Range("H77").Value = 2362
$ws.Range("I77").Value = 1470.4286
$ws.Range("K77").Value = 7352.143
$ws.Range("M77").Value = -2984.143
$ws.Range("H116").Value = 513.75
$ws.Range("I116").Value = 513.75
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 513.75
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1780.25
$ws.Range("N116").ClearContents()
$ws.Range("H121").Value = 27695.73
$ws.Range("J121").Value = 27695.73
$ws.Range("L121").Value = 27695.73
$ws.Range("N121").Value = -31189.73
$ws.Range("H132").Value = 2930.1875
$ws.Range("I132").Value = 1607.7273
$ws.Range("J132").Value = 5839.6
$ws.Range("K132").Value = 4823.1819
$ws.Range("L132").Value = 17518.8
$ws.Range("M132").Value = -2293.1819
$ws.Range("N132").Value = -22578.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 513.75
$ws.Range("I3").Value = 513.75
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 513.75
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -399.75
$ws.Range("N3").ClearContents()
$ws.Range("H107").Value = 2035.6666
$ws.Range("I107").Value = 1930.5
$ws.Range("J107").Value = 2119.8
$ws.Range("K107").Value = 1930.5
$ws.Range("L107").Value = 2119.8
$ws.Range("M107").Value = -10.5
$ws.Range("N107").Value = -5959.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5464.0835
$ws.Range("I31").Value = 2175
$ws.Range("J31").Value = 10068.8
$ws.Range("K31").Value = 2175
$ws.Range("L31").Value = 10068.8
$ws.Range("M31").Value = -1880
$ws.Range("N31").Value = -10658.8
$ws.Range("H34").Value = 5464.0835
$ws.Range("I34").Value = 2175
$ws.Range("J34").Value = 10068.8
$ws.Range("K34").Value = 2175
$ws.Range("L34").Value = 10068.8
$ws.Range("M34").Value = -1973
$ws.Range("N34").Value = -10472.8
$ws.Range("H52").Value = 32633.334
$ws.Range("J52").Value = 32633.334
$ws.Range("L52").Value = 32633.334
$ws.Range("N52").Value = -33221.334
$ws.Range("H81").Value = 30333
$ws.Range("J81").Value = 30333
$ws.Range("L81").Value = 30333
$ws.Range("N81").Value = -32329
$ws.Range("H84").Value = 30333
$ws.Range("J84").Value = 30333
$ws.Range("L84").Value = 90999
$ws.Range("N84").Value = -100983
$ws.Range("H109").Value = 38599
$ws.Range("J109").Value = 38599
$ws.Range("L109").Value = 38599
$ws.Range("N109").Value = -40679

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 219.71428
$ws.Range("I98").Value = 246
$ws.Range("J98").Value = 200
$ws.Range("K98").Value = 738
$ws.Range("L98").Value = 600
$ws.Range("M98").Value = 760
$ws.Range("N98").Value = -3596

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H124").Value = 41880
$ws.Range("J124").Value = 41880
$ws.Range("L124").Value = 41880
$ws.Range("N124").Value = -51700
$ws.Range("H126").Value = 3386.33
$ws.Range("I126").Value = 2789.4932
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 8368.479599999999
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -5898.479599999999
$ws.Range("N126").Value = -19940

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5222.231
$ws.Range("I7").Value = 2048.1667
$ws.Range("J7").Value = 7942.857
$ws.Range("K7").Value = 2048.1667
$ws.Range("L7").Value = 7942.857
$ws.Range("M7").Value = -1936.1667
$ws.Range("N7").Value = -8166.857
$ws.Range("H40").Value = 6687.9473
$ws.Range("I40").Value = 5768.8823
$ws.Range("J40").Value = 14500
$ws.Range("K40").Value = 5768.8823
$ws.Range("L40").Value = 14500
$ws.Range("M40").Value = -5632.8823
$ws.Range("N40").Value = -14772
$ws.Range("H46").Value = 1540.0834
$ws.Range("J46").Value = 2511.2
$ws.Range("L46").Value = 2511.2
$ws.Range("N46").Value = -2887.2
$ws.Range("H126").Value = 5222.231
$ws.Range("I126").Value = 2048.1667
$ws.Range("J126").Value = 7942.857
$ws.Range("K126").Value = 6144.500100000001
$ws.Range("L126").Value = 23828.571
$ws.Range("M126").Value = -3674.500100000001
$ws.Range("N126").Value = -28768.571
$ws.Range("H132").Value = 3661.5715
$ws.Range("I132").Value = 2805.4482
$ws.Range("K132").Value = 8416.3446
$ws.Range("M132").Value = -5886.3446

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 63661.285
$ws.Range("J46").Value = 63661.285
$ws.Range("L46").Value = 63661.285
$ws.Range("N46").Value = -64123.285
$ws.Range("H81").Value = 1999.875
$ws.Range("I81").Value = 1999.875
$ws.Range("K81").Value = 3999.75
$ws.Range("M81").Value = -2938.75
$ws.Range("H84").Value = 1999.875
$ws.Range("I84").Value = 1999.875
$ws.Range("K84").Value = 19998.75
$ws.Range("M84").Value = -14694.75
$ws.Range("H122").Value = 4042
$ws.Range("I122").Value = 2948.4375
$ws.Range("K122").Value = 8845.3125
$ws.Range("M122").Value = -6395.3125
$ws.Range("H123").Value = 38951.75
$ws.Range("J123").Value = 38951.75
$ws.Range("L123").Value = 38951.75
$ws.Range("N123").Value = -48751.75
$ws.Range("H125").Value = 39841.668
$ws.Range("J125").Value = 39841.668
$ws.Range("L125").Value = 39841.668
$ws.Range("N125").Value = -49681.668
$ws.Range("H126").Value = 2239.2593
$ws.Range("I126").Value = 1003.9
$ws.Range("J126").Value = 5768.857
$ws.Range("K126").Value = 3011.7
$ws.Range("L126").Value = 17306.571
$ws.Range("M126").Value = -541.6999999999998
$ws.Range("N126").Value = -22246.571
$ws.Range("H134").Value = 63661.285
$ws.Range("J134").Value = 63661.285
$ws.Range("L134").Value = 190983.855
$ws.Range("N134").Value = -196053.855

